$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 322.1
$ws.Range("I43").Value = 263.5
$ws.Range("J43").Value = 347.2143
$ws.Range("K43").Value = 263.5
$ws.Range("L43").Value = 347.2143
$ws.Range("M43").Value = -194.5
$ws.Range("N43").Value = -485.2143

$ws.Range("H76").Value = 3047.3684
$ws.Range("I76").Value = 3064.2856
$ws.Range("J76").Value = 3000
$ws.Range("K76").Value = 3064.2856
$ws.Range("L76").Value = 3000
$ws.Range("M76").Value = -2749.2856
$ws.Range("N76").Value = -3630

$ws.Range("H79").Value = 3047.3684
$ws.Range("I79").Value = 3064.2856
$ws.Range("J79").Value = 3000
$ws.Range("K79").Value = 3064.2856
$ws.Range("L79").Value = 3000
$ws.Range("M79").Value = -1972.2856
$ws.Range("N79").Value = -5184

$ws.Range("H98").Value = 891.1923
$ws.Range("I98").Value = 931.34784
$ws.Range("J98").Value = 583.3333
$ws.Range("K98").Value = 931.34784
$ws.Range("L98").Value = 583.3333
$ws.Range("M98").Value = 566.65216
$ws.Range("N98").Value = -3579.3333

$ws.Range("H116").Value = 4000
$ws.Range("I116").Value = 5000
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 5000
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = -1558
$ws.Range("N116").Value = -9884

$ws.Range("H122").Value = 891.1923
$ws.Range("I122").Value = 931.34784
$ws.Range("J122").Value = 583.3333
$ws.Range("K122").Value = 2794.04352
$ws.Range("L122").Value = 1749.9999
$ws.Range("M122").Value = -344.0435200000002
$ws.Range("N122").Value = -6649.9999

$ws.Range("H137").Value = 2421.6
$ws.Range("I137").Value = 2660.9092
$ws.Range("J137").Value = 2233.5715
$ws.Range("K137").Value = 7982.7276
$ws.Range("L137").Value = 6700.7145
$ws.Range("M137").Value = -5432.7276
$ws.Range("N137").Value = -11800.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 2632.3333
$ws.Range("I3").Value = 1266.6666
$ws.Range("J3").Value = 3998
$ws.Range("K3").Value = 1266.6666
$ws.Range("L3").Value = 3998
$ws.Range("M3").Value = -1151.6666
$ws.Range("N3").Value = -4228

$ws.Range("H32").Value = 7500.061
$ws.Range("I32").Value = 3723.1538
$ws.Range("J32").Value = 21941.176
$ws.Range("K32").Value = 3723.1538
$ws.Range("L32").Value = 21941.176
$ws.Range("M32").Value = -3436.1538
$ws.Range("N32").Value = -22515.176

$ws.Range("H61").Value = 1747.8572
$ws.Range("I61").Value = 1333.75
$ws.Range("J61").Value = 2300
$ws.Range("K61").Value = 1333.75
$ws.Range("L61").Value = 2300
$ws.Range("M61").Value = -1121.75
$ws.Range("N61").Value = -2724

$ws.Range("H92").Value = 22068
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 22068
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 22068
$ws.Range("N92").Value = -27060

$ws.Range("H132").Value = 2717.3928
$ws.Range("I132").Value = 2405.3333
$ws.Range("J132").Value = 2865.2104
$ws.Range("K132").Value = 7215.999899999999
$ws.Range("L132").Value = 8595.6312
$ws.Range("M132").Value = -4685.999899999999
$ws.Range("N132").Value = -13655.6312

$ws.Range("H136").Value = 1747.8572
$ws.Range("I136").Value = 1333.75
$ws.Range("J136").Value = 2300
$ws.Range("K136").Value = 4001.25
$ws.Range("L136").Value = 6900
$ws.Range("M136").Value = -1451.25
$ws.Range("N136").Value = -12000

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 3157
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 3157
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 3157
$ws.Range("N2").Value = -3383

$ws.Range("H4").Value = 3220.8667
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 3220.8667
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 3220.8667
$ws.Range("N4").Value = -3444.8667

$ws.Range("H5").Value = 371.5625
$ws.Range("I5").Value = 352
$ws.Range("J5").Value = 380.45456
$ws.Range("K5").Value = 352
$ws.Range("L5").Value = 380.45456
$ws.Range("M5").Value = -240
$ws.Range("N5").Value = -604.45456

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 3426.4707
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 3426.4707
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 10279.4121
$ws.Range("N100").Value = -11901.4121

$ws.Range("H113").Value = 1378061.9
$ws.Range("I113").Value = 4329831.5
$ws.Range("J113").Value = 569.2
$ws.Range("K113").Value = 12989494.5
$ws.Range("L113").Value = 1707.6
$ws.Range("M113").Value = -12987324.5
$ws.Range("N113").Value = -6047.6

$ws.Range("H126").Value = 3993.3635
$ws.Range("I126").Value = 1991.25
$ws.Range("J126").Value = 9332.333000000001
$ws.Range("K126").Value = 5973.75
$ws.Range("L126").Value = 27996.999
$ws.Range("M126").Value = -1033.75
$ws.Range("N126").Value = -37876.999

$ws.Range("H132").Value = 2375.6428
$ws.Range("I132").Value = 2688
$ws.Range("J132").Value = 2015.2307
$ws.Range("K132").Value = 24192
$ws.Range("L132").Value = 18137.0763
$ws.Range("M132").Value = -21662
$ws.Range("N132").Value = -23197.0763

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 2537.2
$ws.Range("I4").Value = 1563.3334
$ws.Range("J4").Value = 3998
$ws.Range("K4").Value = 1563.3334
$ws.Range("L4").Value = 3998
$ws.Range("M4").Value = -1451.3334
$ws.Range("N4").Value = -4222

$ws.Range("H5").Value = 11950
$ws.Range("I5").Value = 1250
$ws.Range("J5").Value = 14625
$ws.Range("K5").Value = 1250
$ws.Range("L5").Value = 14625
$ws.Range("M5").Value = -1138
$ws.Range("N5").Value = -14849

$ws.Range("H7").Value = 2200
$ws.Range("I7").Value = 1600
$ws.Range("J7").Value = 4000
$ws.Range("K7").Value = 1600
$ws.Range("L7").Value = 4000
$ws.Range("M7").Value = -1488
$ws.Range("N7").Value = -4224

$ws.Range("H8").Value = 2200
$ws.Range("I8").Value = 1600
$ws.Range("J8").Value = 4000
$ws.Range("K8").Value = 1600
$ws.Range("L8").Value = 4000
$ws.Range("M8").Value = -1461
$ws.Range("N8").Value = -4278

$ws.Range("H122").Value = 1707.6666
$ws.Range("I122").Value = 1578.6154
$ws.Range("J122").Value = 2043.2
$ws.Range("K122").Value = 4735.8462
$ws.Range("L122").Value = 6129.6
$ws.Range("M122").Value = -2285.8462
$ws.Range("N122").Value = -11029.6

$ws.Range("H132").Value = 2661.463
$ws.Range("I132").Value = 2858.1614
$ws.Range("J132").Value = 2396.348
$ws.Range("K132").Value = 8574.484199999999
$ws.Range("L132").Value = 7189.044
$ws.Range("M132").Value = -6044.484199999999
$ws.Range("N132").Value = -12249.044

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 168833.33
$ws.Range("I2").Value = 483333.34
$ws.Range("J2").Value = 11583.333
$ws.Range("K2").Value = 483333.34
$ws.Range("L2").Value = 11583.333
$ws.Range("M2").Value = -483221.34
$ws.Range("N2").Value = -11807.333

$ws.Range("H40").Value = 1807.5883
$ws.Range("I40").Value = 1774.4546
$ws.Range("J40").Value = 1868.3334
$ws.Range("K40").Value = 1774.4546
$ws.Range("L40").Value = 1868.3334
$ws.Range("M40").Value = -1638.4546
$ws.Range("N40").Value = -2140.3334

$ws.Range("H122").Value = 3316.5
$ws.Range("I122").Value = 3702.4167
$ws.Range("J122").Value = 2985.7144
$ws.Range("K122").Value = 11107.2501
$ws.Range("L122").Value = 8957.143199999999
$ws.Range("M122").Value = -8657.250100000001
$ws.Range("N122").Value = -13857.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M11").ClearContents()
$ws.Range("H11").Value = 30000
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 30000
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 30000
$ws.Range("N11").Value = -30284
